$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# E..T numeric columns.
$rows = @(
    @("ECs",  "Ptn", "Sdc3", "ECs",   1, 0.3333333333333333, 0.06520933333333333, 0.195628, 0.007116089623850998, 0.007116089623850999, 3, 1, 28.19948866666667, 84.598466, 0.7357427920402423, 0.7357427920402422, 1.838869856294222, 16.549828706648, 0.005235611648260731, 0.00523561164826073),
    @("ECs",  "Ptn", "Sdc3", "FAPs",  1, 0.3333333333333333, 0.06520933333333333, 0.195628, 0.007116089623850998, 0.007116089623850999, 3, 1, 7.778025666666667, 23.334077, 0.2029336910395279, 0.2029336910395278, 0.5071998683728889, 4.564798815356, 0.001444094333136168, 0.001444094333136168),
    @("ECs",  "Ptn", "Sdc3", "MuSCs", 1, 0.3333333333333333, 0.06520933333333333, 0.195628, 0.007116089623850998, 0.007116089623850999, 3, 1, 2.350402666666667, 7.051208000000001, 0.0613235169202299, 0.06132351692022989, 0.1532681909582222, 1.379413718624, 0.0004363836424540991, 0.000436383642454099),
    @("FAPs", "Ptn", "Sdc3", "ECs",   3, 1, 3.270036666666666, 9.81011, 0.3568488252184601, 0.3568488252184602, 3, 1, 28.19948866666667, 84.598466, 0.7357427920402423, 0.7357427920402422, 92.21336192125111, 829.92025729126, 0.2625489510025102, 0.2625489510025103),
    @("FAPs", "Ptn", "Sdc3", "FAPs",  3, 1, 3.270036666666666, 9.81011, 0.3568488252184601, 0.3568488252184602, 3, 1, 7.778025666666667, 23.334077, 0.2029336910395279, 0.2029336910395278, 25.43442912427444, 228.90986211847, 0.07241664924470145, 0.07241664924470147),
    @("FAPs", "Ptn", "Sdc3", "MuSCs", 3, 1, 3.270036666666666, 9.81011, 0.3568488252184601, 0.3568488252184602, 3, 1, 2.350402666666667, 7.051208000000001, 0.0613235169202299, 0.06132351692022989, 7.685902901431112, 69.17312611288001, 0.0218832249712484, 0.0218832249712484),
    @("MuSCs","Ptn", "Sdc3", "ECs",   3, 1, 5.828401, 17.485203, 0.6360350851576888, 0.6360350851576889, 3, 1, 28.19948866666667, 84.598466, 0.7357427920402423, 0.7357427920402422, 164.3579279442887, 1479.221351498598, 0.4679582293894712, 0.4679582293894712),
    @("MuSCs","Ptn", "Sdc3", "FAPs",  3, 1, 5.828401, 17.485203, 0.6360350851576888, 0.6360350851576889, 3, 1, 7.778025666666667, 23.334077, 0.2029336910395279, 0.2029336910395278, 45.33345257362567, 408.001073162631, 0.1290729474616902, 0.1290729474616902),
    @("MuSCs","Ptn", "Sdc3", "MuSCs", 3, 1, 5.828401, 17.485203, 0.6360350851576888, 0.6360350851576889, 3, 1, 2.350402666666667, 7.051208000000001, 0.0613235169202299, 0.06132351692022989, 13.69908925280267, 123.291803275224, 0.0390039083065274, 0.03900390830652739)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
